$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 12894.786
$ws.Range("I76").Value = 20264.834
$ws.Range("J76").Value = 7367.25
$ws.Range("K76").Value = 20264.834
$ws.Range("L76").Value = 7367.25
$ws.Range("M76").Value = -19949.834
$ws.Range("N76").Value = -7997.25
$ws.Range("H79").Value = 12894.786
$ws.Range("I79").Value = 20264.834
$ws.Range("J79").Value = 7367.25
$ws.Range("K79").Value = 20264.834
$ws.Range("L79").Value = 7367.25
$ws.Range("M79").Value = -19172.834
$ws.Range("N79").Value = -9551.25
$ws.Range("H92").Value = 1409.25
$ws.Range("I92").Value = 1373.3334
$ws.Range("K92").Value = 1373.3334
$ws.Range("M92").Value = -125.3334
$ws.Range("H135").Value = 5814477
$ws.Range("I135").Value = 447.60526
$ws.Range("K135").Value = 4028.44734
$ws.Range("M135").Value = -1493.44734
$ws.Range("H137").Value = 11389.619
$ws.Range("I137").Value = 5249.278
$ws.Range("J137").Value = 48231.668
$ws.Range("K137").Value = 15747.834
$ws.Range("L137").Value = 144695.004
$ws.Range("M137").Value = -13197.834
$ws.Range("N137").Value = -149795.004
$ws.Range("H138").Value = 3655.2307
$ws.Range("J138").Value = 12142
$ws.Range("L138").Value = 36426
$ws.Range("N138").Value = -46706
$ws.Range("H141").Value = 958.67346
$ws.Range("I141").Value = 883.587
$ws.Range("K141").Value = 2650.761
$ws.Range("M141").Value = 2529.239

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3574.6938
$ws.Range("I32").Value = 2816.8635
$ws.Range("K32").Value = 2816.8635
$ws.Range("M32").Value = -2529.8635
$ws.Range("H74").Value = 35755404
$ws.Range("I74").Value = 43527964
$ws.Range("K74").Value = 43527964
$ws.Range("M74").Value = -43527090
$ws.Range("H77").Value = 35755404
$ws.Range("I77").Value = 43527964
$ws.Range("K77").Value = 217639820
$ws.Range("M77").Value = -217635452
$ws.Range("H110").Value = 1819.6923
$ws.Range("I110").Value = 1320.7
$ws.Range("J110").Value = 3483
$ws.Range("K110").Value = 1320.7
$ws.Range("L110").Value = 3483
$ws.Range("M110").Value = 724.3
$ws.Range("N110").Value = -7573

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 19519.709
$ws.Range("J86").Value = 53262.875
$ws.Range("L86").Value = 53262.875
$ws.Range("N86").Value = -55508.875
$ws.Range("H89").Value = 19519.709
$ws.Range("J89").Value = 53262.875
$ws.Range("L89").Value = 266314.375
$ws.Range("N89").Value = -277546.375
$ws.Range("H94").Value = 2821.7646
$ws.Range("I94").Value = 2638.5715
$ws.Range("J94").Value = 2950
$ws.Range("K94").Value = 2638.5715
$ws.Range("L94").Value = 2950
$ws.Range("M94").Value = -2187.5715
$ws.Range("N94").Value = -3852
$ws.Range("H105").Value = 10730.167
$ws.Range("I105").Value = 15278.286
$ws.Range("J105").Value = 4362.8
$ws.Range("K105").Value = 15278.286
$ws.Range("L105").Value = 4362.8
$ws.Range("M105").Value = -13531.286
$ws.Range("N105").Value = -7856.8

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38466084
$ws.Range("J31").Value = 76929700
$ws.Range("L31").Value = 76929700
$ws.Range("N31").Value = -76930290
$ws.Range("H34").Value = 38466084
$ws.Range("J34").Value = 76929700
$ws.Range("L34").Value = 76929700
$ws.Range("N34").Value = -76930104
$ws.Range("H99").Value = 5699.4644
$ws.Range("I99").Value = 5890.421
$ws.Range("J99").Value = 5296.3335
$ws.Range("K99").Value = 5890.421
$ws.Range("L99").Value = 5296.3335
$ws.Range("M99").Value = -4392.421
$ws.Range("N99").Value = -8292.333500000001
$ws.Range("H107").Value = 4444
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 4444
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 4444
$ws.Range("N107").Value = -8284
$ws.Range("H126").Value = 5699.4644
$ws.Range("I126").Value = 5890.421
$ws.Range("J126").Value = 5296.3335
$ws.Range("K126").Value = 17671.263
$ws.Range("L126").Value = 15889.0005
$ws.Range("M126").Value = -15201.263
$ws.Range("N126").Value = -20829.0005
$ws.Range("H132").Value = 143582.4
$ws.Range("I132").Value = 194775.95
$ws.Range("J132").Value = 9199.375
$ws.Range("K132").Value = 584327.8500000001
$ws.Range("L132").Value = 27598.125
$ws.Range("M132").Value = -581797.8500000001
$ws.Range("N132").Value = -32658.125
$ws.Range("H134").Value = 1813.3125
$ws.Range("I134").Value = 1525.25
$ws.Range("J134").Value = 2677.5
$ws.Range("K134").Value = 4575.75
$ws.Range("L134").Value = 8032.5
$ws.Range("M134").Value = -2040.75
$ws.Range("N134").Value = -13102.5
$ws.Range("M107").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 977.7727
$ws.Range("J107").Value = 1148.8
$ws.Range("L107").Value = 3446.4
$ws.Range("N107").Value = -7286.4
$ws.Range("H133").Value = 12474.6
$ws.Range("H134").Value = 1856.1765
$ws.Range("I134").Value = 2115.875
$ws.Range("J134").Value = 1232.9
$ws.Range("K134").Value = 6347.625
$ws.Range("L134").Value = 3698.7
$ws.Range("M134").Value = -1277.625
$ws.Range("N134").Value = -13838.7
$ws.Range("H139").Value = 2143
$ws.Range("I139").Value = 2031.75
$ws.Range("K139").Value = 6095.25
$ws.Range("M139").Value = -955.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 11543273
$ws.Range("I126").Value = 5887366.5
$ws.Range("K126").Value = 17662099.5
$ws.Range("M126").Value = -17659629.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2626.9714
$ws.Range("I7").Value = 2144.7334
$ws.Range("K7").Value = 2144.7334
$ws.Range("M7").Value = -2032.7334
$ws.Range("H61").Value = 9334.666999999999
$ws.Range("I61").Value = 8004
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 8004
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -7802
$ws.Range("N61").Value = -10404
$ws.Range("H109").Value = 67000
$ws.Range("J109").Value = 60000
$ws.Range("L109").Value = 60000
$ws.Range("N109").Value = -62774
$ws.Range("H113").Value = 9334.666999999999
$ws.Range("I113").Value = 8004
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 8004
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -5834
$ws.Range("N113").Value = -14340
$ws.Range("H126").Value = 2626.9714
$ws.Range("I126").Value = 2144.7334
$ws.Range("K126").Value = 6434.2002
$ws.Range("M126").Value = -3964.2002
$ws.Range("H132").Value = 50011390
$ws.Range("I132").Value = 7930.8667
$ws.Range("J132").Value = 200021780
$ws.Range("K132").Value = 23792.6001
$ws.Range("L132").Value = 600065340
$ws.Range("M132").Value = -21262.6001
$ws.Range("N132").Value = -600070400
$ws.Range("H136").Value = 2268.276
$ws.Range("I136").Value = 1941.8334
$ws.Range("K136").Value = 5825.5002
$ws.Range("M136").Value = -3275.5002

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 979.36365
$ws.Range("I81").Value = 998.7778
$ws.Range("K81").Value = 1997.5556
$ws.Range("M81").Value = -936.5555999999999
$ws.Range("H84").Value = 979.36365
$ws.Range("I84").Value = 998.7778
$ws.Range("K84").Value = 9987.778
$ws.Range("M84").Value = -4683.778
$ws.Range("H107").Value = 228.66667
$ws.Range("I107").Value = 120.2
$ws.Range("J107").Value = 306.14285
$ws.Range("K107").Value = 360.6
$ws.Range("L107").Value = 918.4285500000001
$ws.Range("M107").Value = 1559.4
$ws.Range("N107").Value = -4758.428550000001
$ws.Range("H113").Value = 906.75
$ws.Range("I113").Value = 787
$ws.Range("K113").Value = 2361
$ws.Range("M113").Value = -191
$ws.Range("H117").Value = 60000
$ws.Range("J117").Value = 60000
$ws.Range("L117").Value = 60000
$ws.Range("N117").Value = -69178
$ws.Range("H132").Value = 1506.525
$ws.Range("I132").Value = 1454.2368
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 4362.7104
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -1832.7104
$ws.Range("N132").Value = -12560
$ws.Range("H136").Value = 2264.8235
$ws.Range("I136").Value = 781.6667
$ws.Range("K136").Value = 2345.0001
$ws.Range("M136").Value = 204.9998999999998
